$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 355 (August Red / Venus, week of 2022-02-21),
# shifting all rows from 355 onward down to 359 onward. The 4 new rows carry a new
# weekly price report (2022-12-23, date serial 44918) for two new varieties:
# "Artic Pride" (Especial/Primera) and "Early John" (Especial/Primera).
$ws.Rows("355:358").Insert()

# Row 355: Artic Pride / Especial
$ws.Range("A355").Value = 8
$ws.Range("B355").Value = "Terminal La Palmera de La Serena"
$ws.Range("C355").Value = "Coquimbo"
$ws.Range("D355").Value = 44918
$ws.Range("E355").Value = 4
$ws.Range("F355").Value = "Fruta"
$ws.Range("G355").Value = 100103
$ws.Range("H355").Value = "Frutos de hueso (carozo)"
$ws.Range("I355").Value = 100103006
$ws.Range("J355").Value = "Nectarín"
$ws.Range("K355").Value = "Artic Pride"
$ws.Range("L355").Value = "Especial"
$ws.Range("M355").Value = 16
$ws.Range("N355").Value = 300000
$ws.Range("O355").Value = 310000
$ws.Range("P355").Value = 305000
$ws.Range("Q355").Value = "$/bins (420 kilos)"
$ws.Range("R355").Value = "Región de O'Higgins"
$ws.Range("S355").Value = 726
$ws.Range("T355").Value = 420

# Row 356: Artic Pride / Primera
$ws.Range("A356").Value = 8
$ws.Range("B356").Value = "Terminal La Palmera de La Serena"
$ws.Range("C356").Value = "Coquimbo"
$ws.Range("D356").Value = 44918
$ws.Range("E356").Value = 4
$ws.Range("F356").Value = "Fruta"
$ws.Range("G356").Value = 100103
$ws.Range("H356").Value = "Frutos de hueso (carozo)"
$ws.Range("I356").Value = 100103006
$ws.Range("J356").Value = "Nectarín"
$ws.Range("K356").Value = "Artic Pride"
$ws.Range("L356").Value = "Primera"
$ws.Range("M356").Value = 16
$ws.Range("N356").Value = 240000
$ws.Range("O356").Value = 250000
$ws.Range("P356").Value = 245000
$ws.Range("Q356").Value = "$/bins (420 kilos)"
$ws.Range("R356").Value = "Región de O'Higgins"
$ws.Range("S356").Value = 583
$ws.Range("T356").Value = 420

# Row 357: Early John / Especial
$ws.Range("A357").Value = 8
$ws.Range("B357").Value = "Terminal La Palmera de La Serena"
$ws.Range("C357").Value = "Coquimbo"
$ws.Range("D357").Value = 44918
$ws.Range("E357").Value = 4
$ws.Range("F357").Value = "Fruta"
$ws.Range("G357").Value = 100103
$ws.Range("H357").Value = "Frutos de hueso (carozo)"
$ws.Range("I357").Value = 100103006
$ws.Range("J357").Value = "Nectarín"
$ws.Range("K357").Value = "Early John"
$ws.Range("L357").Value = "Especial"
$ws.Range("M357").Value = 10
$ws.Range("N357").Value = 500000
$ws.Range("O357").Value = 510000
$ws.Range("P357").Value = 505000
$ws.Range("Q357").Value = "$/bins (420 kilos)"
$ws.Range("R357").Value = "Región de O'Higgins"
$ws.Range("S357").Value = 1202
$ws.Range("T357").Value = 420

# Row 358: Early John / Primera
$ws.Range("A358").Value = 8
$ws.Range("B358").Value = "Terminal La Palmera de La Serena"
$ws.Range("C358").Value = "Coquimbo"
$ws.Range("D358").Value = 44918
$ws.Range("E358").Value = 4
$ws.Range("F358").Value = "Fruta"
$ws.Range("G358").Value = 100103
$ws.Range("H358").Value = "Frutos de hueso (carozo)"
$ws.Range("I358").Value = 100103006
$ws.Range("J358").Value = "Nectarín"
$ws.Range("K358").Value = "Early John"
$ws.Range("L358").Value = "Primera"
$ws.Range("M358").Value = 20
$ws.Range("N358").Value = 470000
$ws.Range("O358").Value = 480000
$ws.Range("P358").Value = 475000
$ws.Range("Q358").Value = "$/bins (420 kilos)"
$ws.Range("R358").Value = "Región de O'Higgins"
$ws.Range("S358").Value = 1131
$ws.Range("T358").Value = 420
